$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append " order By ..." suffix to the query texts in B2, B3, B4
$ws.Range("B2").Value = $ws.Range("B2").Value2 + "`n order By ss.study_subject_id ASC LIMIT 100 "
$ws.Range("B3").Value = $ws.Range("B3").Value2 + "`n order By samp.sample_id ASC LIMIT 100"
$ws.Range("B4").Value = $ws.Range("B4").Value2 + "`n order By f.file_name ASC LIMIT 100"

# Row heights grow because the wrapped text got longer (Excel autofit result)
$ws.Rows.Item(2).RowHeight = 345.6
$ws.Rows.Item(3).RowHeight = 259.2
$ws.Rows.Item(4).RowHeight = 244.8

# Update selection / view to match after-state (A4 topLeftCell, C4 active selection)
$ws.Range("C4").Select()
